$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item("Tabla1")
$col = $tbl.ListColumns.Add()

$ws.Range("E8").Value = "DONE"
$ws.Range("E18").Value = "DONE"
$ws.Range("E31").Value = "DONE"
$ws.Range("E1").Value = "Column1"

$ws.Range("E11").Select() | Out-Null
